$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "template"
$wb.Worksheets.Item("Sheet1").Name = "template"

$ws = $wb.Worksheets.Item("expected")

# --- table.update(): shrink tables by moving them further down the sheet.
# Read every old table's values up front (their old/new footprints can
# overlap), then clear_contents the old footprints (instead of deleting
# cells, which would shift unrelated rows), then resize the tables and
# write the relocated values into their new footprints.

$lo7 = $ws.ListObjects("Table715")
$oldRange7 = $ws.Range("A24:E28")
$vals7 = $oldRange7.Value2

$lo8 = $ws.ListObjects("Table1116")
$oldRange8 = $ws.Range("A37:D41")
$vals8 = $oldRange8.Value2

$oldRange7.ClearContents()
$oldRange8.ClearContents()

$lo7.Resize($ws.Range("A39:E43"))
$lo8.Resize($ws.Range("A52:D56"))

$ws.Range("A39:E43").Value2 = $vals7
$ws.Range("A52:D56").Value2 = $vals8

# Restore the selection on the "expected" sheet
$ws.Activate()
$ws.Range("A23").Select()
